# Insert a new data row at row 330 (pushes the existing rows 330-440 down
# to 331-441, which is exactly the shift seen throughout the diff: every
# row's "after" content equals the row above's "before" content, and a
# brand-new row appears at the top of the shifted block with a new last
# row 441 at the bottom carrying what used to be row 440's data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("330").Insert()

$ws.Range("A330").Value = 10
$ws.Range("B330").Value = "Vega Modelo de Temuco"
$ws.Range("C330").Value = "La Araucanía"
$ws.Range("D330").Value = 44468
$ws.Range("E330").Value = 9
$ws.Range("F330").Value = 100114001
$ws.Range("G330").Value = "Papa"
$ws.Range("H330").Value = "Asterix"
$ws.Range("I330").Value = "1a (guarda)"
$ws.Range("J330").Value = 300
$ws.Range("K330").Value = 9000
$ws.Range("L330").Value = 9000
$ws.Range("M330").Value = 9000
$ws.Range("N330").Value = "`$/malla 25 kilos"
$ws.Range("O330").Value = "Provincia de Cautín"
$ws.Range("P330").Value = 360
$ws.Range("Q330").Value = 25
$ws.Range("R330").Value = "Hortaliza"
